# Apply the "Update Structurizer, Utilizer, main" edit:
#  - re-purpose/rename header columns Y..AD
#  - flip a couple of boolean flags in row 2
#  - replace several boolean cells with free-text reasoning / markers
#  - blank out a couple of cells that are no longer used

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 1: header renames ----
$ws.Range("Y1").Value  = "verdict"
$ws.Range("Z1").Value  = "reason"
$ws.Range("AA1").Value = "媒體影響_value"
$ws.Range("AB1").Value = "媒體影響_type"
$ws.Range("AC1").Value = "量刑爭議_value"
$ws.Range("AD1").Value = "量刑爭議_type"

# ---- Row 2 ----
$ws.Range("L2").Value = $true
$ws.Range("T2").Value = $true

$ws.Range("Z2").Value = "。根據裁定書的理由，法院認為本案存在不適宜行國民參與審判的情事，包括被告已對被訴事實作有罪陳述、案件情節可能對被害人家屬造成二度傷害，以及案件涉及的證據可能對國民法官造成過大刺激，符合國民法官法第6條第1項第4款及第5款的規定。"
$ws.Range("AA2").Value = $true
$ws.Range("AB2").Value = "bool"
$ws.Range("AC2").Value = ""
$ws.Range("AD2").Value = ""

# ---- Row 3 ----
$ws.Range("Z3").Value = "。根據裁定書的理由，雖然本案符合國民法官法第5條第1項第2款的適用條件，但被告已認罪且與被害人家屬達成和解，且檢辯雙方對於量刑無重大爭執，因此法院認為不行國民參與審判為適當。"
$ws.Range("AA3").Value = ""
$ws.Range("AB3").Value = ""
$ws.Range("AD3").Value = "bool"
